# "Generate Report for Archive"
# The localization-status report was regenerated: the row describing
# f267f0d1-0300-46ae-b972-c5a06a9f0ff6 and the row describing
# 633c5ece-5139-4489-b415-aa0b99d7bbec swapped places (rows 6 and 7 on
# every sheet), and f267f0d1's status moved from "Ready for handoff" to
# "In Translation" now that it sits in row 6.

$wb = $excel.ActiveWorkbook

function Set-RowSixAndSeven($ws, $a6, $b6, $c6, $d6, $a7, $b7, $c7, $d7, $hasD) {
    $ws.Range("A6").Value = $a6
    $ws.Range("B6").Value = $b6
    $ws.Range("C6").Value = $c6
    $ws.Range("A7").Value = $a7
    $ws.Range("B7").Value = $b7
    $ws.Range("C7").Value = $c7

    if ($hasD) {
        $ws.Range("D6").Value = $d6
        $ws.Range("D7").Value = $d7
    }

    # Hyperlinks stay bound to their cell (and keep pointing at the same
    # external target), only the displayed text follows the new content.
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq '$A$6') {
            $h.TextToDisplay = $a6
        } elseif ($addr -eq '$A$7') {
            $h.TextToDisplay = $a7
        } elseif ($hasD -and $addr -eq '$C$6') {
            $h.TextToDisplay = $c6
        } elseif ($hasD -and $addr -eq '$C$7') {
            $h.TextToDisplay = $c7
        }
    }
}

# --- Sheet 1: Overview ---
$wsOverview = $wb.Worksheets.Item(1)
Set-RowSixAndSeven $wsOverview `
    "f267f0d1-0300-46ae-b972-c5a06a9f0ff6.md" "In Translation" "In Translation" "" `
    "633c5ece-5139-4489-b415-aa0b99d7bbec.md" "Ready for handoff" "Ready for handoff" "" `
    $false

# --- Sheet 2: zh-cn ---
$wsZhCn = $wb.Worksheets.Item(2)
Set-RowSixAndSeven $wsZhCn `
    "f267f0d1-0300-46ae-b972-c5a06a9f0ff6.md" "In Translation" `
    "f267f0d1-0300-46ae-b972-c5a06a9f0ff6.0001c6190457a4bc7d05ec8578fa22b2ddb4258c.zh-cn.xlf" "2016-02-24 06:49:51" `
    "633c5ece-5139-4489-b415-aa0b99d7bbec.md" "Ready for handoff" `
    "633c5ece-5139-4489-b415-aa0b99d7bbec.348d93e84a68523c1d12fe2ba726f3c1d928c2c4.zh-cn.xlf" "2016-02-24 06:44:29" `
    $true

# --- Sheet 3: de-de ---
$wsDeDe = $wb.Worksheets.Item(3)
Set-RowSixAndSeven $wsDeDe `
    "f267f0d1-0300-46ae-b972-c5a06a9f0ff6.md" "In Translation" `
    "f267f0d1-0300-46ae-b972-c5a06a9f0ff6.0001c6190457a4bc7d05ec8578fa22b2ddb4258c.de-de.xlf" "2016-02-24 06:50:03" `
    "633c5ece-5139-4489-b415-aa0b99d7bbec.md" "Ready for handoff" `
    "633c5ece-5139-4489-b415-aa0b99d7bbec.348d93e84a68523c1d12fe2ba726f3c1d928c2c4.de-de.xlf" "2016-02-24 06:44:41" `
    $true
